$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New volunteer-hours entry logged: Start Time / Finish Time / Minutes Logged
$ws.Cells.Item(22, 1).Value = "10:58PM 1-19-2018"
$ws.Cells.Item(22, 2).Value = "11:55PM 1-19-2018"
$ws.Cells.Item(22, 3).Value = 57

# Move the active selection to reflect where the user ended up editing
[void]$ws.Range("C20").Select()
